# 🔄 Actualización automática del tracker
# Fill in results ("resultado") and profit ("profit") for rows that were
# previously pending, and normalize two event_id cells (A89, A90) that were
# stored as text to be proper numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => (resultado, profit)
$updates = @{
    69 = @("Fallo",   -1)
    75 = @("Acierto",  0.62)
    76 = @("Fallo",   -1)
    77 = @("Fallo",   -1)
    83 = @("Acierto",  2.4)
    85 = @("Fallo",   -1)
    87 = @("Acierto",  1.5)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 7).Value = $vals[0]
    $ws.Cells.Item($row, 8).Value = $vals[1]
}

# Normalize event_id cells that were stored as text to numeric values.
$ws.Cells.Item(89, 1).Value = 14265609
$ws.Cells.Item(90, 1).Value = 14359058
